$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the ranking table from 8 entries (rows 2-9) to 11 entries (rows 2-12) ---
# New leaderboard order/values (rank, placement text, nick name, points):
#  1  1º  bona final boss   9999   (unchanged)
#  2  2º  bona              2461
#  3  3º  felipe            2446
#  4  4º  teste             2436
#  5  5º  nan               2430
#  6  6º  matheus           2430
#  7  7º  Diegowl           2344
#  8  8º  Esther linda      2310
#  9  9º  last dance        1977
# 10 10º  oi                1922
# 11 11º  felca             1422

# First, clone the bold/border/centered style used by A2:A9 down into the
# three brand-new rank cells (A10:A12) by copying formats from A9.
$ws.Range("A9").Copy()
$ws.Range("A10:A12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$rankData = @(
  @{ Row = 2;  Rank = 1;  Place = "1º";  Nick = "bona final boss"; Pts = "9999" },
  @{ Row = 3;  Rank = 2;  Place = "2º";  Nick = "bona";             Pts = "2461" },
  @{ Row = 4;  Rank = 3;  Place = "3º";  Nick = "felipe";           Pts = "2446" },
  @{ Row = 5;  Rank = 4;  Place = "4º";  Nick = "teste";            Pts = "2436" },
  @{ Row = 6;  Rank = 5;  Place = "5º";  Nick = "nan";              Pts = "2430" },
  @{ Row = 7;  Rank = 6;  Place = "6º";  Nick = "matheus";          Pts = "2430" },
  @{ Row = 8;  Rank = 7;  Place = "7º";  Nick = "Diegowl";          Pts = "2344" },
  @{ Row = 9;  Rank = 8;  Place = "8º";  Nick = "Esther linda";     Pts = "2310" },
  @{ Row = 10; Rank = 9;  Place = "9º";  Nick = "last dance";       Pts = "1977" },
  @{ Row = 11; Rank = 10; Place = "10º"; Nick = "oi";                Pts = "1922" },
  @{ Row = 12; Rank = 11; Place = "11º"; Nick = "felca";             Pts = "1422" }
)

foreach ($entry in $rankData) {
    $r = $entry.Row

    $ws.Cells.Item($r, 1).Value = $entry.Rank
    $ws.Cells.Item($r, 2).Value = $entry.Place
    $ws.Cells.Item($r, 3).Value = $entry.Nick

    # Points are stored as text (not numbers) in the sheet, so force the
    # quote-prefix text entry, then strip the resulting "quote prefix" style
    # back to Normal so no stray numeric formatting is left behind.
    $ws.Cells.Item($r, 4).Value = "'" + $entry.Pts
    $ws.Cells.Item($r, 4).Style = "Normal"
}
